# Remove `pax` from databases
# Rename header labels in both worksheets of USE_TYPE_PROPERTIES.xlsx:
#   Occ_m2pax  -> Occ_m2p
#   Vww_lpdpax -> Vww_ldp
#   Vw_lpdpax  -> Vw_ldp
#   X_ghpax    -> X_ghp
#   Qs_Wpax    -> Qs_Wp
#   Ve_lpspax  -> Ve_lsp

$wb = $excel.ActiveWorkbook

$wsLoads = $wb.Worksheets.Item("INTERNAL_LOADS")
$wsComfort = $wb.Worksheets.Item("INDOOR_COMFORT")

# INTERNAL_LOADS sheet header row (row 1)
$wsLoads.Range("B1").Value = "Occ_m2p"
$wsLoads.Range("C1").Value = "Qs_Wp"
$wsLoads.Range("D1").Value = "X_ghp"
$wsLoads.Range("I1").Value = "Vww_ldp"
$wsLoads.Range("J1").Value = "Vw_ldp"

# INDOOR_COMFORT sheet header row (row 1)
$wsComfort.Range("F1").Value = "Ve_lsp"
